# Update the "K" column (column G) values in the save_data sheet for
# littell_zack.xlsx. Per the commit message, the save_data pipeline was
# regenerated to compute K (strikeouts) differently ("use K instead of
# Strike#"), so the raw per-game K values in column G (rows 2-34) are
# replaced with the newly calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$newK = [ordered]@{
    2  = 7
    3  = 7
    4  = 5
    5  = 4
    6  = 2
    7  = 3
    8  = 5
    9  = 3
    10 = 5
    11 = 6
    12 = 3
    13 = 5
    14 = 4
    15 = 6
    16 = 5
    17 = 2
    18 = 3
    19 = 2
    20 = 9
    21 = 7
    22 = 3
    23 = 3
    24 = 6
    25 = 6
    26 = 7
    27 = 7
    28 = 4
    29 = 5
    30 = 6
    31 = 3
    32 = 6
    33 = 1
    34 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
